$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '27.886.14'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.630.75'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.19'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0880'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").Value = '1.861.64'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").Value = '1.621.35'
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.561'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '27.884.14'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.48'
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.11'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  -3.88%  '
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.00'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.90'
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.51'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0482'
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.41'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").Value = '1.393.73'
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("E36").Value = '  +9.08%  '
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0171'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.558'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.80'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("D46").Value = '1.771.23'
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.46'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  +0.06%  '
